# "Update US model copy" -- refresh the Max Biofuel Blends assumptions and
# move the active selection around to match the author's final editing
# session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Data edit on the "Max Biofuel Blends" sheet: the biodiesel/biogas ramp
#    now starts in 2020 (was 2021), and the max biogas blend by 2050 is
#    raised from 0.5 to 0.6. All FORECAST.LINEAR / TREND formulas on this
#    sheet, plus every per-fuel sheet that pulls from row 14/17 of this
#    sheet, recalculate automatically.
# ---------------------------------------------------------------------------
$wsMBB = $wb.Worksheets.Item("Max Biofuel Blends")
$wsMBB.Range("A6").Value = 2020
$wsMBB.Range("A9").Value = 2020
$wsMBB.Range("B10").Value = 0.6

# ---------------------------------------------------------------------------
# 2. Leave a selection behind on "MPoEFUbVT-ships-frgt-dslveh" (it was the
#    second-to-last sheet the author touched) before moving on.
# ---------------------------------------------------------------------------
$wsShips = $wb.Worksheets.Item("MPoEFUbVT-ships-frgt-dslveh")
$wsShips.Activate()
$wsShips.Range("H5").Select()

# ---------------------------------------------------------------------------
# 3. Finish on the "Max Biofuel Blends" sheet with A10 selected -- this is
#    the sheet/cell that ends up active (and tabSelected) when the file is
#    saved.
# ---------------------------------------------------------------------------
$wsMBB.Activate()
$wsMBB.Range("A10").Select()
